$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New quantities ("disponible" stock count) for column F, rows 2-19.
# These replace the placeholder text value "1" that was shared by every row.
$values = [ordered]@{
    2  = 9
    3  = 10
    4  = 30
    5  = 6
    6  = 2
    7  = 6
    8  = 23
    9  = 2
    10 = 7
    11 = 16
    12 = 3
    13 = 5
    14 = 0
    15 = 15
    16 = 12
    17 = 8
    18 = 5
    19 = 8
}

# Apply the new number format (integer "0") and alignment (centered
# horizontally, top vertically) to the first cell only. Doing this on a
# single cell lets the engine reuse/mutate one style record instead of
# minting a fresh one for every property change.
$first = $ws.Range("F2")
$first.NumberFormat = "0"
$first.HorizontalAlignment = -4108   # xlCenter
$first.VerticalAlignment = -4160     # xlTop

# Propagate that exact style to the rest of the column in a single
# format-only paste so no duplicate style records are created.
$first.Copy()
$ws.Range("F3:F19").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Now write the actual numeric values.
foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}

# Update the paper size / orientation used for printing.
$ps = $ws.PageSetup
$ps.PaperSize = 9   # xlPaperA4
$ps.Orientation = 1 # xlPortrait

# Reflect the author's final selection.
$ws.Range("F19").Select()

$wb.Save()
